# DemoTransactions.xlsx update
# - tweak a few existing debit amounts
# - remove a couple of stray "Sub category" (F) values that shouldn't be there
# - append new transaction rows for May/June/July 2022
# - update the sheet's scroll/selection state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Tweak existing debit amounts
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = 136.61
$ws.Cells.Item(13, 4).Value = 227.67
$ws.Cells.Item(23, 4).Value = 126.54
$ws.Cells.Item(29, 4).Value = 54.99

# ---------------------------------------------------------------------------
# 2. Remove stray F values on rows 3 and 6
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 6).ClearContents()
$ws.Cells.Item(6, 6).ClearContents()

# ---------------------------------------------------------------------------
# 3. Append new transaction rows (45-77)
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=45; Date=44711; B='AMAZON UK'; C='Books from Amazon'; D=35.99; F=2; G=$null; Style=1 },
    @{ Row=46; Date=44683; B='Hotel at Booking.com'; C='Hotel in Madrid'; D=216.35; F=9; G=901; Style=1 },
    @{ Row=47; Date=44684; B='CO-OP GROUP'; C='Groceries'; D=23.35; F=1; G=$null; Style=1 },
    @{ Row=48; Date=44685; B='WESTMINSTER COUNCIL'; C='Council tax'; D=158; F=5; G=$null; Style=1 },
    @{ Row=49; Date=44686; B='TFL TRAVELCARD'; C='Monthly railcard'; D=120; F=6; G=601; Style=1 },
    @{ Row=50; Date=44689; B='CO-OP GROUP'; C='Groceries'; D=32.32; F=1; G=$null; Style=1 },
    @{ Row=51; Date=44693; B='TESCO GROUP'; C='Groceries'; D=65.25; F=1; G=$null; Style=1 },
    @{ Row=52; Date=44697; B='M&A Plumbing and Heating'; C='Boiler repair'; D=85; F=12; G=1201; Style=1 },
    @{ Row=53; Date=44700; B='EDF ENERGY'; C='Electricity bills'; D=23.36; F=4; G=402; Style=1 },
    @{ Row=54; Date=44701; B='THAMES WATER'; C='Water bills'; D=18.25; F=4; G=401; Style=1 },
    @{ Row=55; Date=44704; B='Cineworld'; C='Cinema'; D=24.99; F=10; G=$null; Style=1 },
    @{ Row=56; Date=44705; B='WAITROSE & PARTNERS'; C='Groceries'; D=25.14; F=1; G=$null; Style=1 },
    @{ Row=57; Date=44715; B='CO-OP GROUP'; C='Groceries'; D=35.36; F=1; G=$null; Style=1 },
    @{ Row=58; Date=44716; B='WESTMINSTER COUNCIL'; C='Council tax'; D=158; F=5; G=$null; Style=1 },
    @{ Row=59; Date=44717; B='TFL TRAVELCARD'; C='Monthly railcard'; D=125; F=6; G=601; Style=1 },
    @{ Row=60; Date=44720; B='CO-OP GROUP'; C='Groceries'; D=37.35; F=1; G=$null; Style=1 },
    @{ Row=61; Date=44724; B='TESCO GROUP'; C='Groceries'; D=59.35; F=1; G=$null; Style=1 },
    @{ Row=62; Date=44728; B='M+L Electricity'; C='Cooker repair'; D=50; F=12; G=1201; Style=1 },
    @{ Row=63; Date=44731; B='EDF ENERGY'; C='Electricity bills'; D=35.24; F=4; G=402; Style=1 },
    @{ Row=64; Date=44732; B='THAMES WATER'; C='Water bills'; D=25.36; F=4; G=401; Style=1 },
    @{ Row=65; Date=44735; B='ARCHERY'; C='Archery Club Membership'; D=70; F=10; G=$null; Style=1 },
    @{ Row=66; Date=44736; B='WAITROSE & PARTNERS'; C='Groceries'; D=38.35; F=1; G=$null; Style=1 },
    @{ Row=67; Date=44742; B='JUST EAT'; C='Takeaway'; D=45; F=8; G=$null; Style=2 },
    @{ Row=68; Date=44746; B='WESTMINSTER COUNCIL'; C='Council tax'; D=130; F=5; G=$null; Style=1 },
    @{ Row=69; Date=44750; B='CO-OP GROUP'; C='Groceries'; D=25.36; F=1; G=$null; Style=1 },
    @{ Row=70; Date=44754; B='TESCO GROUP'; C='Groceries'; D=45.87; F=1; G=$null; Style=1 },
    @{ Row=71; Date=44758; B='M+L Electricity'; C='Light bulb'; D=15; F=12; G=1201; Style=1 },
    @{ Row=72; Date=44761; B='EDF ENERGY'; C='Electricity bills'; D=24.14; F=4; G=402; Style=1 },
    @{ Row=73; Date=44762; B='THAMES WATER'; C='Water bills'; D=37.25; F=4; G=401; Style=1 },
    @{ Row=74; Date=44765; B='Cineworld'; C='Cinema'; D=12.99; F=10; G=$null; Style=1 },
    @{ Row=75; Date=44766; B='WAITROSE & PARTNERS'; C='Groceries'; D=24.32; F=1; G=$null; Style=1 },
    @{ Row=76; Date=44772; B='AMAZON UK'; C='Fan from Amazon'; D=79.99; F=2; G=$null; Style=1 },
    @{ Row=77; Date=44772; B='JUST EAT'; C='Takeaway'; D=15; F=8; G=$null; Style=2 }
)

# Template rows that already carry the right number formats / styles:
#  - row 34 -> style 1 (plain date format), columns A:G
#  - row 22 -> style 2 (right aligned date format), columns A:G
foreach ($item in $newRows) {
    $r = $item.Row
    if ($item.Style -eq 2) {
        [void]$ws.Range("A22:G22").Copy()
    } else {
        [void]$ws.Range("A34:G34").Copy()
    }
    [void]$ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = $item.F
    if ($item.G -ne $null) {
        $ws.Cells.Item($r, 7).Value = $item.G
    } else {
        $ws.Cells.Item($r, 7).ClearContents()
    }
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Update view state - scroll to row 34 and select D53
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D53").Select()
